$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-29 Friday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-11-30 Saturday", 2) | Out-Null
$d.Content.Find.Execute("89-5=", $true, $true, $false, $false, $false, $true, 1, $false, "90+9=", 2) | Out-Null
$d.Content.Find.Execute("77-73=", $true, $true, $false, $false, $false, $true, 1, $false, "25+29=", 2) | Out-Null
$d.Content.Find.Execute("78-57=", $true, $true, $false, $false, $false, $true, 1, $false, "20+45=", 2) | Out-Null
$d.Content.Find.Execute("89-34=", $true, $true, $false, $false, $false, $true, 1, $false, "0+41=", 2) | Out-Null
$d.Content.Find.Execute("61+22=", $true, $true, $false, $false, $false, $true, 1, $false, "5+8=", 2) | Out-Null
$d.Content.Find.Execute("2+94=", $true, $true, $false, $false, $false, $true, 1, $false, "96-84=", 2) | Out-Null
$d.Content.Find.Execute("88-73=", $true, $true, $false, $false, $false, $true, 1, $false, "42+37=", 2) | Out-Null
$d.Content.Find.Execute("18+81=", $true, $true, $false, $false, $false, $true, 1, $false, "68+10=", 2) | Out-Null
$d.Content.Find.Execute("88-19=", $true, $true, $false, $false, $false, $true, 1, $false, "59-21=", 2) | Out-Null
$d.Content.Find.Execute("96-82=", $true, $true, $false, $false, $false, $true, 1, $false, "36+37=", 2) | Out-Null
$d.Content.Find.Execute("87-74=", $true, $true, $false, $false, $false, $true, 1, $false, "37+6=", 2) | Out-Null
$d.Content.Find.Execute("92-69=", $true, $true, $false, $false, $false, $true, 1, $false, "50+24=", 2) | Out-Null
$d.Content.Find.Execute("76-38=", $true, $true, $false, $false, $false, $true, 1, $false, "30+43=", 2) | Out-Null
$d.Content.Find.Execute("58+20=", $true, $true, $false, $false, $false, $true, 1, $false, "30+9=", 2) | Out-Null
$d.Content.Find.Execute("17+39=", $true, $true, $false, $false, $false, $true, 1, $false, "7+59=", 2) | Out-Null
$d.Content.Find.Execute("73-30=", $true, $true, $false, $false, $false, $true, 1, $false, "95-65=", 2) | Out-Null
$d.Content.Find.Execute("76-5=", $true, $true, $false, $false, $false, $true, 1, $false, "93-85=", 2) | Out-Null
$d.Content.Find.Execute("43+35=", $true, $true, $false, $false, $false, $true, 1, $false, "54-24=", 2) | Out-Null
$d.Content.Find.Execute("24+39=", $true, $true, $false, $false, $false, $true, 1, $false, "40+37=", 2) | Out-Null
$d.Content.Find.Execute("32+50=", $true, $true, $false, $false, $false, $true, 1, $false, "61+13=", 2) | Out-Null
$d.Content.Find.Execute("22-12=", $true, $true, $false, $false, $false, $true, 1, $false, "43-11=", 2) | Out-Null
$d.Content.Find.Execute("70+28=", $true, $true, $false, $false, $false, $true, 1, $false, "15+8=", 2) | Out-Null
$d.Content.Find.Execute("39+11=", $true, $true, $false, $false, $false, $true, 1, $false, "20+19=", 2) | Out-Null
$d.Content.Find.Execute("69-23=", $true, $true, $false, $false, $false, $true, 1, $false, "79-21=", 2) | Out-Null
$d.Content.Find.Execute("43+56=", $true, $true, $false, $false, $false, $true, 1, $false, "17-2=", 2) | Out-Null
$d.Content.Find.Execute("82+7=", $true, $true, $false, $false, $false, $true, 1, $false, "75-5=", 2) | Out-Null
$d.Content.Find.Execute("49+31=", $true, $true, $false, $false, $false, $true, 1, $false, "63-41=", 2) | Out-Null
$d.Content.Find.Execute("13+27=", $true, $true, $false, $false, $false, $true, 1, $false, "94-29=", 2) | Out-Null
$d.Content.Find.Execute("78-43=", $true, $true, $false, $false, $false, $true, 1, $false, "61+36=", 2) | Out-Null
$d.Content.Find.Execute("7+45=", $true, $true, $false, $false, $false, $true, 1, $false, "21-4=", 2) | Out-Null
$d.Content.Find.Execute("59-56=", $true, $true, $false, $false, $false, $true, 1, $false, "53+2=", 2) | Out-Null
$d.Content.Find.Execute("81-29=", $true, $true, $false, $false, $false, $true, 1, $false, "79-0=", 2) | Out-Null
$d.Content.Find.Execute("84-57=", $true, $true, $false, $false, $false, $true, 1, $false, "11+50=", 2) | Out-Null
$d.Content.Find.Execute("54+9=", $true, $true, $false, $false, $false, $true, 1, $false, "56-19=", 2) | Out-Null
$d.Content.Find.Execute("58+23=", $true, $true, $false, $false, $false, $true, 1, $false, "95-88=", 2) | Out-Null
$d.Content.Find.Execute("73-63=", $true, $true, $false, $false, $false, $true, 1, $false, "81-35=", 2) | Out-Null
$d.Content.Find.Execute("65-40=", $true, $true, $false, $false, $false, $true, 1, $false, "58-26=", 2) | Out-Null
$d.Content.Find.Execute("36+18=", $true, $true, $false, $false, $false, $true, 1, $false, "50+33=", 2) | Out-Null
$d.Content.Find.Execute("78-59=", $true, $true, $false, $false, $false, $true, 1, $false, "1+92=", 2) | Out-Null
$d.Content.Find.Execute("62-13=", $true, $true, $false, $false, $false, $true, 1, $false, "85+0=", 2) | Out-Null
$d.Content.Find.Execute("16+42=", $true, $true, $false, $false, $false, $true, 1, $false, "11+40=", 2) | Out-Null
$d.Content.Find.Execute("67-52=", $true, $true, $false, $false, $false, $true, 1, $false, "55+34=", 2) | Out-Null
$d.Content.Find.Execute("70+25=", $true, $true, $false, $false, $false, $true, 1, $false, "43+37=", 2) | Out-Null
$d.Content.Find.Execute("71-63=", $true, $true, $false, $false, $false, $true, 1, $false, "17+9=", 2) | Out-Null
$d.Content.Find.Execute("18+0=", $true, $true, $false, $false, $false, $true, 1, $false, "10+8=", 2) | Out-Null
$d.Content.Find.Execute("46-31=", $true, $true, $false, $false, $false, $true, 1, $false, "20+13=", 2) | Out-Null
$d.Content.Find.Execute("38+42=", $true, $true, $false, $false, $false, $true, 1, $false, "20+40=", 2) | Out-Null
$d.Content.Find.Execute("52+43=", $true, $true, $false, $false, $false, $true, 1, $false, "51-13=", 2) | Out-Null
$d.Content.Find.Execute("29-0=", $true, $true, $false, $false, $false, $true, 1, $false, "60-30=", 2) | Out-Null
$d.Content.Find.Execute("89-58=", $true, $true, $false, $false, $false, $true, 1, $false, "13-0=", 2) | Out-Null
$d.Content.Find.Execute("56-0=", $true, $true, $false, $false, $false, $true, 1, $false, "70-7=", 2) | Out-Null
$d.Content.Find.Execute("98-31=", $true, $true, $false, $false, $false, $true, 1, $false, "38+33=", 2) | Out-Null
$d.Content.Find.Execute("42+28=", $true, $true, $false, $false, $false, $true, 1, $false, "69-13=", 2) | Out-Null
$d.Content.Find.Execute("10+23=", $true, $true, $false, $false, $false, $true, 1, $false, "90-61=", 2) | Out-Null
$d.Content.Find.Execute("57-32=", $true, $true, $false, $false, $false, $true, 1, $false, "36+21=", 2) | Out-Null
$d.Content.Find.Execute("84-56=", $true, $true, $false, $false, $false, $true, 1, $false, "7+41=", 2) | Out-Null
$d.Content.Find.Execute("24+59=", $true, $true, $false, $false, $false, $true, 1, $false, "3+93=", 2) | Out-Null
$d.Content.Find.Execute("60+39=", $true, $true, $false, $false, $false, $true, 1, $false, "90-41=", 2) | Out-Null
$d.Content.Find.Execute("90-19=", $true, $true, $false, $false, $false, $true, 1, $false, "75+6=", 2) | Out-Null
$d.Content.Find.Execute("40+28=", $true, $true, $false, $false, $false, $true, 1, $false, "2+80=", 2) | Out-Null
$d.Content.Find.Execute("49-27=", $true, $true, $false, $false, $false, $true, 1, $false, "76-75=", 2) | Out-Null
$d.Content.Find.Execute("62-52=", $true, $true, $false, $false, $false, $true, 1, $false, "76-59=", 2) | Out-Null
$d.Content.Find.Execute("17+37=", $true, $true, $false, $false, $false, $true, 1, $false, "80-45=", 2) | Out-Null
$d.Content.Find.Execute("81-34=", $true, $true, $false, $false, $false, $true, 1, $false, "13+18=", 2) | Out-Null
$d.Content.Find.Execute("24+50=", $true, $true, $false, $false, $false, $true, 1, $false, "42+15=", 2) | Out-Null
$d.Content.Find.Execute("14+41=", $true, $true, $false, $false, $false, $true, 1, $false, "53-14=", 2) | Out-Null
$d.Content.Find.Execute("85-82=", $true, $true, $false, $false, $false, $true, 1, $false, "92-86=", 2) | Out-Null
$d.Content.Find.Execute("10+67=", $true, $true, $false, $false, $false, $true, 1, $false, "97-61=", 2) | Out-Null
$d.Content.Find.Execute("22-16=", $true, $true, $false, $false, $false, $true, 1, $false, "40+35=", 2) | Out-Null
$d.Content.Find.Execute("24-22=", $true, $true, $false, $false, $false, $true, 1, $false, "28+32=", 2) | Out-Null
$d.Content.Find.Execute("86-47=", $true, $true, $false, $false, $false, $true, 1, $false, "47+44=", 2) | Out-Null
$d.Content.Find.Execute("29+31=", $true, $true, $false, $false, $false, $true, 1, $false, "46-37=", 2) | Out-Null
$d.Content.Find.Execute("60-50=", $true, $true, $false, $false, $false, $true, 1, $false, "62-18=", 2) | Out-Null
$d.Content.Find.Execute("44-10=", $true, $true, $false, $false, $false, $true, 1, $false, "40-39=", 2) | Out-Null
$d.Content.Find.Execute("59-42=", $true, $true, $false, $false, $false, $true, 1, $false, "31+42=", 2) | Out-Null
$d.Content.Find.Execute("33+28=", $true, $true, $false, $false, $false, $true, 1, $false, "65-42=", 2) | Out-Null
$d.Content.Find.Execute("63-27=", $true, $true, $false, $false, $false, $true, 1, $false, "5+79=", 2) | Out-Null
$d.Content.Find.Execute("22+16=", $true, $true, $false, $false, $false, $true, 1, $false, "11+18=", 2) | Out-Null
$d.Content.Find.Execute("87-38=", $true, $true, $false, $false, $false, $true, 1, $false, "30+45=", 2) | Out-Null
$d.Content.Find.Execute("86-51=", $true, $true, $false, $false, $false, $true, 1, $false, "47+49=", 2) | Out-Null
$d.Content.Find.Execute("86-37=", $true, $true, $false, $false, $false, $true, 1, $false, "68-44=", 2) | Out-Null
$d.Content.Find.Execute("77-19=", $true, $true, $false, $false, $false, $true, 1, $false, "46+32=", 2) | Out-Null
$d.Content.Find.Execute("81+8=", $true, $true, $false, $false, $false, $true, 1, $false, "65-15=", 2) | Out-Null
$d.Content.Find.Execute("77-42=", $true, $true, $false, $false, $false, $true, 1, $false, "73-53=", 2) | Out-Null
$d.Content.Find.Execute("86-32=", $true, $true, $false, $false, $false, $true, 1, $false, "38-14=", 2) | Out-Null
$d.Content.Find.Execute("14-0=", $true, $true, $false, $false, $false, $true, 1, $false, "23+46=", 2) | Out-Null
$d.Content.Find.Execute("28-6=", $true, $true, $false, $false, $false, $true, 1, $false, "21+24=", 2) | Out-Null
$d.Content.Find.Execute("64-42=", $true, $true, $false, $false, $false, $true, 1, $false, "26+52=", 2) | Out-Null
$d.Content.Find.Execute("97-89=", $true, $true, $false, $false, $false, $true, 1, $false, "19+63=", 2) | Out-Null
$d.Content.Find.Execute("96-37=", $true, $true, $false, $false, $false, $true, 1, $false, "31+43=", 2) | Out-Null
$d.Content.Find.Execute("88-2=", $true, $true, $false, $false, $false, $true, 1, $false, "50+48=", 2) | Out-Null
$d.Content.Find.Execute("0+24=", $true, $true, $false, $false, $false, $true, 1, $false, "69-24=", 2) | Out-Null
$d.Content.Find.Execute("19+62=", $true, $true, $false, $false, $false, $true, 1, $false, "30-1=", 2) | Out-Null
$d.Content.Find.Execute("62-15=", $true, $true, $false, $false, $false, $true, 1, $false, "17+22=", 2) | Out-Null
$d.Content.Find.Execute("13+85=", $true, $true, $false, $false, $false, $true, 1, $false, "49+18=", 2) | Out-Null
$d.Content.Find.Execute("74-24=", $true, $true, $false, $false, $false, $true, 1, $false, "7+15=", 2) | Out-Null
$d.Content.Find.Execute("13+84=", $true, $true, $false, $false, $false, $true, 1, $false, "23+61=", 2) | Out-Null
$d.Content.Find.Execute("42-8=", $true, $true, $false, $false, $false, $true, 1, $false, "44-3=", 2) | Out-Null
$d.Content.Find.Execute("63-45=", $true, $true, $false, $false, $false, $true, 1, $false, "75+21=", 2) | Out-Null
$d.Content.Find.Execute("32-11=", $true, $true, $false, $false, $false, $true, 1, $false, "62-20=", 2) | Out-Null
